$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 7458188
$ws.Range("I116").Value = 8587805
$ws.Range("J116").Value = 2719.8
$ws.Range("K116").Value = 8587805
$ws.Range("L116").Value = 2719.8
$ws.Range("M116").Value = -8584363
$ws.Range("N116").Value = -9603.799999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 3032.8572
$ws.Range("I102").Value = 3667.5
$ws.Range("J102").Value = 2186.6667
$ws.Range("K102").Value = 3667.5
$ws.Range("L102").Value = 2186.6667
$ws.Range("M102").Value = -2045.5
$ws.Range("N102").Value = -5430.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 22729546
$ws.Range("I16").Value = 27779810
$ws.Range("J16").Value = 3356.5
$ws.Range("K16").Value = 27779810
$ws.Range("L16").Value = 3356.5
$ws.Range("M16").Value = -27779523
$ws.Range("N16").Value = -3930.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3793740.2
$ws.Range("I31").Value = 5584288.5
$ws.Range("J31").Value = 1991.4706
$ws.Range("K31").Value = 5584288.5
$ws.Range("L31").Value = 1991.4706
$ws.Range("M31").Value = -5583993.5
$ws.Range("N31").Value = -2581.4706

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3793740.2
$ws.Range("I34").Value = 5584288.5
$ws.Range("J34").Value = 1991.4706
$ws.Range("K34").Value = 5584288.5
$ws.Range("L34").Value = 1991.4706
$ws.Range("M34").Value = -5584086.5
$ws.Range("N34").Value = -2395.4706

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 22729546
$ws.Range("I113").Value = 27779810
$ws.Range("J113").Value = 3356.5
$ws.Range("K113").Value = 27779810
$ws.Range("L113").Value = 3356.5
$ws.Range("M113").Value = -27777640
$ws.Range("N113").Value = -7696.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 22222888
$ws.Range("I122").Value = 22222888
$ws.Range("K122").Value = 66668664
$ws.Range("M122").Value = -66666214

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 9260022
$ws.Range("I131").Value = 288.4
$ws.Range("K131").Value = 865.1999999999999
$ws.Range("M131").Value = 4174.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2631.5789
$ws.Range("I126").Value = 1609.091
$ws.Range("J126").Value = 4037.5
$ws.Range("K126").Value = 4827.272999999999
$ws.Range("L126").Value = 12112.5
$ws.Range("M126").Value = -2357.272999999999
$ws.Range("N126").Value = -17052.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H127").Value = 16331.5
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 16331.5
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 16331.5
$ws.Range("N127").Value = -26251.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 0

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H129").Value = 34222
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 34222
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 34222
$ws.Range("N129").Value = -44222

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H130").Value = 80780
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 80780
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 80780
$ws.Range("N130").Value = -90820

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H131").Value = 62750
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 62750
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 62750
$ws.Range("N131").Value = -72830

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2983.025
$ws.Range("I132").Value = 2854.3928
$ws.Range("J132").Value = 3283.1667
$ws.Range("K132").Value = 8563.178400000001
$ws.Range("L132").Value = 9849.500100000001
$ws.Range("M132").Value = -6033.178400000001
$ws.Range("N132").Value = -14909.5001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 25000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 25000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 25000
$ws.Range("N133").Value = -35120

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 57465.2
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 57465.2
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 172395.6
$ws.Range("N134").Value = -177465.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 12500
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 12500
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 12500
$ws.Range("N135").Value = -22640

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 30809.334
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 30809.334
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 92428.00199999999
$ws.Range("N136").Value = -97528.00199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H139").Value = 31000
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 31000
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 31000
$ws.Range("N139").Value = -41280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H140").Value = 49766.668
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 49766.668
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 49766.668
$ws.Range("N140").Value = -60126.668

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H141").Value = 32871.5
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 32871.5
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 32871.5
$ws.Range("N141").Value = -43231.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 953.6129
$ws.Range("I16").Value = 536.96155
$ws.Range("J16").Value = 3120.2
$ws.Range("K16").Value = 536.96155
$ws.Range("L16").Value = 3120.2
$ws.Range("M16").Value = -366.96155
$ws.Range("N16").Value = -3460.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H124").Value = 14942.667
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 14942.667
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 14942.667
$ws.Range("N124").Value = -24762.667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H125").Value = 45450
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 45450
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 45450
$ws.Range("N125").Value = -55290

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H127").Value = 37838.332
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 37838.332
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 37838.332
$ws.Range("N127").Value = -47758.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H128").Value = 29800
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 29800
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 29800
$ws.Range("N128").Value = -39760

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H129").Value = 26500
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 26500
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 26500
$ws.Range("N129").Value = -36500

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H130").Value = 25000
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 25000
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 25000
$ws.Range("N130").Value = -35040

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H131").Value = 65100
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 65100
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 65100
$ws.Range("N131").Value = -75180

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 13304906
$ws.Range("I132").Value = 28421154
$ws.Range("J132").Value = 2608.36
$ws.Range("K132").Value = 85263462
$ws.Range("L132").Value = 7825.08
$ws.Range("M132").Value = -85260932
$ws.Range("N132").Value = -12885.08

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 16608.666
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 16608.666
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 16608.666
$ws.Range("N133").Value = -21668.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H134").Value = 36291.125
$ws.Range("I134").Value = 10000
$ws.Range("J134").Value = 40047
$ws.Range("K134").Value = 10000
$ws.Range("L134").Value = 40047
$ws.Range("M134").Value = -4930
$ws.Range("N134").Value = -50187

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H135").Value = 50409.668
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 50409.668
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 50409.668
$ws.Range("N135").Value = -60549.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4267.15
$ws.Range("I136").Value = 5698.115
$ws.Range("J136").Value = 1609.6428
$ws.Range("K136").Value = 17094.345
$ws.Range("L136").Value = 4828.928400000001
$ws.Range("M136").Value = -14544.345
$ws.Range("N136").Value = -9928.928400000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H137").Value = 58080
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 58080
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 58080
$ws.Range("N137").Value = -68280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H138").Value = 19464.5
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 19464.5
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 19464.5
$ws.Range("N138").Value = -29744.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H139").Value = 28233.334
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 28233.334
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 28233.334
$ws.Range("N139").Value = -38513.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H140").Value = 20000
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 20000
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 20000
$ws.Range("N140").Value = -30360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H141").Value = 48028.75
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 48028.75
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 48028.75
$ws.Range("N141").Value = -58388.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 17101
$ws.Range("J135").Value = 17101
$ws.Range("L135").Value = 17101
$ws.Range("N135").Value = -27241

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H140").Value = 37656.363
$ws.Range("J140").Value = 37656.363
$ws.Range("L140").Value = 37656.363
$ws.Range("N140").Value = -48016.363

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 49724.5
$ws.Range("J141").Value = 49724.5
$ws.Range("L141").Value = 49724.5
$ws.Range("N141").Value = -60084.5
